# Updates the cryptos price/volume data on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number must keep their original
# text representation (they were stored as text in the source data), so we
# force a Text number format before writing the value to stop Excel from
# re-interpreting the string as a numeric value.
$textCells = @("D5", "D8", "D18", "D25", "D26", "D31", "D32", "D33", "D44", "D45", "D46", "D48")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$updates = @(
    @{ Cell = "D2"; Value = "28.260.28" },
    @{ Cell = "E2"; Value = "  +3.75%  " },
    @{ Cell = "D3"; Value = "1.590.25" },
    @{ Cell = "E3"; Value = "  +1.75%  " },
    @{ Cell = "E4"; Value = "  +0.10%  " },
    @{ Cell = "D5"; Value = "213.77" },
    @{ Cell = "E5"; Value = "  +1.41%  " },
    @{ Cell = "E6"; Value = "  +0.94%  " },
    @{ Cell = "E7"; Value = "  +0.05%  " },
    @{ Cell = "D8"; Value = "24.15" },
    @{ Cell = "E8"; Value = "  +8.81%  " },
    @{ Cell = "E9"; Value = "  +0.55%  " },
    @{ Cell = "E10"; Value = "  +0.93%  " },
    @{ Cell = "E11"; Value = "  +2.12%  " },
    @{ Cell = "D12"; Value = "1.818.09" },
    @{ Cell = "E12"; Value = "  +1.78%  " },
    @{ Cell = "D13"; Value = "1.592.92" },
    @{ Cell = "E13"; Value = "  +1.86%  " },
    @{ Cell = "E14"; Value = "  +2.49%  " },
    @{ Cell = "E15"; Value = "  -0.27%  " },
    @{ Cell = "D16"; Value = "28.317.78" },
    @{ Cell = "E16"; Value = "  +3.93%  " },
    @{ Cell = "E17"; Value = "  +2.22%  " },
    @{ Cell = "D18"; Value = "227.87" },
    @{ Cell = "E18"; Value = "  +4.60%  " },
    @{ Cell = "D19"; Value = "0.0₃0709" },
    @{ Cell = "E19"; Value = "  +1.03%  " },
    @{ Cell = "E20"; Value = "  +0.66%  " },
    @{ Cell = "E22"; Value = "  -0.66%  " },
    @{ Cell = "E23"; Value = "  -0.72%  " },
    @{ Cell = "E24"; Value = "  +0.33%  " },
    @{ Cell = "D25"; Value = "151.80" },
    @{ Cell = "E25"; Value = "  +0.28%  " },
    @{ Cell = "D26"; Value = "15.21" },
    @{ Cell = "E26"; Value = "  +1.40%  " },
    @{ Cell = "E27"; Value = "  +0.80%  " },
    @{ Cell = "E28"; Value = "  -0.53%  " },
    @{ Cell = "E29"; Value = "  +0.06%  " },
    @{ Cell = "E30"; Value = "  -0.46%  " },
    @{ Cell = "D31"; Value = "0.0474" },
    @{ Cell = "D32"; Value = "3.24" },
    @{ Cell = "E32"; Value = "  -0.04%  " },
    @{ Cell = "D33"; Value = "3.15" },
    @{ Cell = "E33"; Value = "  -0.87%  " },
    @{ Cell = "D34"; Value = "1.401.43" },
    @{ Cell = "E34"; Value = "  -3.89%  " },
    @{ Cell = "E35"; Value = "  -2.07%  " },
    @{ Cell = "E36"; Value = "  -6.85%  " },
    @{ Cell = "E37"; Value = "  +0.09%  " },
    @{ Cell = "E39"; Value = "  +8.47%  " },
    @{ Cell = "E40"; Value = "  +0.27%  " },
    @{ Cell = "E41"; Value = "  +0.17%  " },
    @{ Cell = "E42"; Value = "  +0.07%  " },
    @{ Cell = "E43"; Value = "  -3.63%  " },
    @{ Cell = "D44"; Value = "1.88" },
    @{ Cell = "E44"; Value = "  +7.08%  " },
    @{ Cell = "D45"; Value = "0.986" },
    @{ Cell = "E45"; Value = "  +1.13%  " },
    @{ Cell = "D46"; Value = "64.30" },
    @{ Cell = "E46"; Value = "  -0.04%  " },
    @{ Cell = "D47"; Value = "1.731.94" },
    @{ Cell = "E47"; Value = "  +1.81%  " },
    @{ Cell = "D48"; Value = "87.57" },
    @{ Cell = "E48"; Value = "  +2.09%  " },
    @{ Cell = "E49"; Value = "  +1.26%  " },
    @{ Cell = "E50"; Value = "  -1.44%  " },
    @{ Cell = "E51"; Value = "  -0.03%  " }
)

foreach ($update in $updates) {
    $ws.Range($update.Cell).Value = $update.Value
}
